# The deck's theme part (ppt/theme/theme1.xml) currently holds the
# "Integral" / "Red Violet" color scheme used by the slide master. The
# target edit replaces its 10 non-black/white slots with the stock
# "Office Theme" palette (dk1/lt1 stay 000000/FFFFFF in both schemes, and
# the fontScheme/fmtScheme are already byte-identical between the two
# theme parts in this deck, so only the clrScheme RGB values need to
# change here).
#
# ThemeColorScheme.Item(n) order, confirmed empirically against the
# package's clrScheme element order:
#   1 dk1  2 lt1  3 dk2  4 lt2  5 accent1 6 accent2 7 accent3
#   8 accent4 9 accent5 10 accent6 11 hlink 12 folHlink

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# index -> target RGB() value (decoded from the target srgbClr hex RRGGBB
# as r | g<<8 | b<<16, matching PowerPoint's VBA RGB() encoding)
$targets = @{
   1 = 0          # dk1      000000
   2 = 16777215   # lt1      FFFFFF
   3 = 6968388    # dk2      44546A
   4 = 15132391   # lt2      E7E6E6
   5 = 13998939   # accent1  5B9BD5
   6 = 3243501    # accent2  ED7D31
   7 = 10855845   # accent3  A5A5A5
   8 = 49407      # accent4  FFC000
   9 = 12874308   # accent5  4472C4
  10 = 4697456    # accent6  70AD47
  11 = 12673797   # hlink    0563C1
  12 = 7491477    # folHlink 954F72
}

for ($i = 1; $i -le $colors.Count; $i++) {
  $colors.Item($i).RGB = $targets[$i]
}
